$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update a few Property/Value pairs ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.1.0"
$meta.Range("B8").Value = "2023-07-10T23:08:03+02:00"
$meta.Range("B10").Value = "No display for ContactDetail"

# --- Sheet "Include from FSIII": remove the concept row that was added ---
$concepts = $wb.Worksheets.Item("Include from FSIII")
$concepts.Rows.Item(13).Delete()
